$d = $word.ActiveDocument

$oldText = "Waktu Kampanye 2018 untuk Perseus: 30 Oktober-8 November dan 29 November-8 Desember"
$newText = "Waktu Kampanye Bootes: 14-23 Mei, 13-22 Juni, 12-21 Juli"

# Every paragraph containing the old campaign-dates sentence (it is split
# across two or three differently-formatted runs, and in one place is
# preceded in the same paragraph by a "www.globeatnight.org" run + line
# break) gets fully collapsed into a single, plain run with the new text
# and no run-level formatting at all, matching the target markup.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$oldText*") {
        $rng = $p.Range
        # Exclude the trailing paragraph mark from the range we clear.
        $rng.End = $rng.End - 1
        # Remove every existing run (and their formatting) in the paragraph.
        $rng.Delete()
        # Insert the replacement as a brand-new, unformatted run.
        $rng.InsertBefore($newText)
    }
}
